$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 'The PowerPoint presentation is exceptional in design and delivery. Slides are cohesive with a consistent design, using appropriate fonts, colors, and objects. The presentation flows smoothly, is engaging, and keeps the audience''s attention. The presenter delivers the content confidently and clearly within the allotted time slot. '
$ws.Range("E5").Value = 'The PowerPoint presentation lacks some cohesiveness in design elements, making it less visually appealing. The delivery is somewhat choppy but generally understandable, and the presentation mostly fits the allotted time. '
$ws.Range("F5").Value = 'The PowerPoint presentation is poorly designed, lacking consistency in fonts, colors, and objects. The delivery is disjointed or difficult to follow, and the presentation significantly exceeds or falls short of the allotted time. '
$ws.Range("G2").Value = 'Missing or no work was submitted.'
$ws.Range("G1").Value = 'Missing'
$ws.Range("C4").Value = 'The visualization stands out in terms of quality, style, clarity, and its ability to help clearly outline how the data science use cases are extracted from the business problem, how they are prioritized and which is the primary use case selected and why. Color, graph choice, labeling, and descriptions are thoughtfully and effectively implemented to drive that argumentation and keep the audience engaged throughout the presentation..'
$ws.Range("E4").Value = 'A visualization was added that was useful to communicate the different use cases to be discussed but the visual could have been improved to better argue and motivate the prioritization of use cases in terms of importance / feasibility. It was not used strategically to emphasize key points in the presentation.'
$ws.Range("F4").Value = 'A visualization is included in the slide deck but only has limited relevance to motivate the selection and prioritization of data science use cases from the business problem at hand. Key elements in the graph like axis labels are missing or hard to read.'
$ws.Range("F3").Value = 'Data science use cases are mentioned but not prioritized or key use cases are missing. It is not clear how some or any of the use cases are relevant for the business problem. A comparison of use case feasibility, resource requirements, timelines and other constraints is completely missing.'
$ws.Range("E3").Value = 'The presentation shows all the relevant data science use cases that could be extracted from business problem. The motivation for the final selection and prioritization of use cases is not clearly fleshed out. The discussion of key considerations like feasibility, resource requirements, timelines, other constraints is not complete.'
$ws.Range("F2").Value = 'It is not clear from the presentation what the final selected use case is. '
$ws.Range("D2").Value = 'The final selected use case is clearly separated from other possible use cases and outlined at the end of the presentation. Verbal and visual presentation are aligned well to communicate the Big Idea / primary use case effectively.'
$ws.Range("C3").Value = 'Visual and verbal presentation align perfectly to elegantly and convincingly walk the audience through the use case selection process, how the use cases were prioritized based on criteria like feasibility, resource requirements, timelines, risks, and/or other constraints and argue effectively why the final use case was selected. The reasoning is persuasive and well designed to also convince non-technical business stakeholders.'
$ws.Range("E2").Value = 'The final selected use case is mentioned but could have been better separated from other use cases mentioned.'
$ws.Range("C2").Value = 'Visually through smart use of text placement, spacing, style, visuals and audibly through speaker''s / speakers'' tone, pace, and other oratory techniques the audience is guided effectively and convincingly to the final selected use case.'
$ws.Range("D4").Value = 'The visualization used was relevant, to the point and added key information to illustrate how the possible data science use cases are motivated by the business problem at and and which use case to prioritize at the end. It used space, color, the choice of graphing style and other elements professionally and effectively to help with the flow of the messaging.'
$ws.Range("D3").Value = 'Visual and verbal presentation components integrate well to make the case for the selection of use cases and which ones should be prioritized and why based on some of the criteria like feasibility, resource requirements, timelines, risks, and/or other constraints. The motivation is easy to follow and is convincing based on the data and ther analysis presented.'
$ws.Range("D5").Value = 'The PowerPoint presentation is well-designed and delivered. Slides exhibit a good degree of cohesion in terms of design elements. The presentation is clear and mostly fluent, and the presenter manages to stay within the allotted time. '

# G3-G5 reuse the same 'Missing or no work was submitted.' text already interned above
$ws.Range("G3").Value = 'Missing or no work was submitted.'
$ws.Range("G4").Value = 'Missing or no work was submitted.'
$ws.Range("G5").Value = 'Missing or no work was submitted.'

# Apply font color (black rgb) to G2:G5 + row 5 (C:G) to get new font/style (fontId=2, xf=4)
$ws.Range("G2:G5").Font.Color = 0
$ws.Range("C5:F5").Font.Color = 0

# Row heights
$ws.Rows.Item(2).RowHeight = 85
$ws.Rows.Item(3).RowHeight = 153
$ws.Rows.Item(4).RowHeight = 153
$ws.Rows.Item(5).RowHeight = 119

# Column G width (26.0 storage units => input 25.17 due to Excel's px-rounding of ColumnWidth)
$ws.Columns.Item(7).ColumnWidth = 25.17

# Selection
$ws.Range("C5").Select()
